$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.397.86"
$ws.Range("E2").Value = "  +1.52%  "
$ws.Range("D3").Value = "1.687.85"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  +0.49%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "218.53"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.97%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.5531"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +8.80%  "
$ws.Range("E7").Value = "  +0.49%  "
$ws.Range("E8").Value = "  +2.36%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06485"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.61%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "22.13"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.74%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07603"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +2.16%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.688.63"
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "4.559"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.22%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.5819"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.09%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.000008474"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.63%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "65.47"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +1.99%  "
$ws.Range("D17").Value = "26.471.85"
$ws.Range("E17").Value = "  +1.42%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "4.949"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.47%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "1.009"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.45%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "10.97"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.88%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "190.98"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.15%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.254"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.31%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "149.77"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +3.70%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.1321"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +10.37%  "
$ws.Range("E26").Value = "  +4.16%  "
$ws.Range("E27").Value = "  +1.36%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.06334"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -3.57%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.397"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +4.48%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.327"
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "3.590"
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.585"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +2.19%  "
$ws.Range("E33").Value = "  +1.32%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.042"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +2.53%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.6250"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +2.15%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.407"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.58%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.721"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.42%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "6.246"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.09%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01636"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +2.61%  "
$ws.Range("D40").Value = "1.116.75"
$ws.Range("E40").Value = "  +2.09%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.8792"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.87%  "
$ws.Range("E42").Value = "  +0.57%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "100.78"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("E44").Value = "  +1.20%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.00000000110"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -3.22%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "57.45"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +2.02%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "8.212"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +2.08%  "
$ws.Range("E48").Value = "  -0.16%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.05283"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.08%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.4301"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.35%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "6.090"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.17%  "
